$wb = $excel.ActiveWorkbook

# ALC!row6
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 1803.4166  # H6: 1815.0834 -> 1803.4166
$ws.Cells.Item(6, 9).Value = 148.2  # I6: 176.2 -> 148.2
$ws.Cells.Item(6, 11).Value = 444.6  # K6: 528.5999999999999 -> 444.6
$ws.Cells.Item(6, 13).Value = -332.6  # M6: -416.5999999999999 -> -332.6

# ALC!row13
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(13, 8).Value = 83339.336  # H13: 52857.43 -> 83339.336
$ws.Cells.Item(13, 10).Value = 83339.336  # J13: 52857.43 -> 83339.336
$ws.Cells.Item(13, 12).Value = 83339.336  # L13: 52857.43 -> 83339.336
$ws.Cells.Item(13, 14).Value = -83677.336  # N13: -53195.43 -> -83677.336

# ALC!row19
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 12987559  # H19: 7143445 -> 12987559
$ws.Cells.Item(19, 9).Value = 35714740  # I19: 11905260 -> 35714740
$ws.Cells.Item(19, 10).Value = 597.1429000000001  # J19: 723.625 -> 597.1429000000001
$ws.Cells.Item(19, 11).Value = 35714740  # K19: 11905260 -> 35714740
$ws.Cells.Item(19, 12).Value = 597.1429000000001  # L19: 723.625 -> 597.1429000000001
$ws.Cells.Item(19, 13).Value = -35714565  # M19: -11905085 -> -35714565
$ws.Cells.Item(19, 14).Value = -947.1429000000001  # N19: -1073.625 -> -947.1429000000001

# ALC!row48
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(48, 8).Value = 2229.75  # H48: 2809.5 -> 2229.75
$ws.Cells.Item(48, 9).Value = 300  # I48: 0 -> 300
$ws.Cells.Item(48, 10).Value = 2873  # J48: 2809.5 -> 2873
$ws.Cells.Item(48, 11).Value = 900  # K48: 0 -> 900
$ws.Cells.Item(48, 12).Value = 8619  # L48: 8428.5 -> 8619
$ws.Cells.Item(48, 13).Value = -608  # M48: None -> -608
$ws.Cells.Item(48, 14).Value = -9203  # N48: -9012.5 -> -9203

# ALC!row56
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(56, 8).Value = 2229.75  # H56: 2809.5 -> 2229.75
$ws.Cells.Item(56, 9).Value = 300  # I56: 0 -> 300
$ws.Cells.Item(56, 10).Value = 2873  # J56: 2809.5 -> 2873
$ws.Cells.Item(56, 11).Value = 900  # K56: 0 -> 900
$ws.Cells.Item(56, 12).Value = 8619  # L56: 8428.5 -> 8619
$ws.Cells.Item(56, 13).Value = -366  # M56: None -> -366
$ws.Cells.Item(56, 14).Value = -9687  # N56: -9496.5 -> -9687

# ALC!row98
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 2791.4075  # H98: 2087.875 -> 2791.4075
$ws.Cells.Item(98, 9).Value = 2760.3809  # I98: 2071.3333 -> 2760.3809
$ws.Cells.Item(98, 10).Value = 2900  # J98: 2137.5 -> 2900
$ws.Cells.Item(98, 11).Value = 2760.3809  # K98: 2071.3333 -> 2760.3809
$ws.Cells.Item(98, 12).Value = 2900  # L98: 2137.5 -> 2900
$ws.Cells.Item(98, 13).Value = -1262.3809  # M98: -573.3332999999998 -> -1262.3809
$ws.Cells.Item(98, 14).Value = -5896  # N98: -5133.5 -> -5896

# ALC!row107
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 669.6429000000001  # H107: 1172.2142 -> 669.6429000000001
$ws.Cells.Item(107, 9).Value = 608.3333  # I107: 1438.125 -> 608.3333
$ws.Cells.Item(107, 10).Value = 780  # J107: 817.6667 -> 780
$ws.Cells.Item(107, 11).Value = 608.3333  # K107: 1438.125 -> 608.3333
$ws.Cells.Item(107, 12).Value = 780  # L107: 817.6667 -> 780
$ws.Cells.Item(107, 13).Value = 1311.6667  # M107: 481.875 -> 1311.6667
$ws.Cells.Item(107, 14).Value = -4620  # N107: -4657.6667 -> -4620

# ALC!row116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 3574.762  # H116: 3327.1738 -> 3574.762
$ws.Cells.Item(116, 9).Value = 3442.8572  # I116: 2890.5 -> 3442.8572
$ws.Cells.Item(116, 10).Value = 3640.7144  # J116: 3663.077 -> 3640.7144
$ws.Cells.Item(116, 11).Value = 3442.8572  # K116: 2890.5 -> 3442.8572
$ws.Cells.Item(116, 12).Value = 3640.7144  # L116: 3663.077 -> 3640.7144
$ws.Cells.Item(116, 13).Value = -0.8571999999999207  # M116: 551.5 -> -0.8571999999999207
$ws.Cells.Item(116, 14).Value = -10524.7144  # N116: -10547.077 -> -10524.7144

# ALC!row122
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 2791.4075  # H122: 2087.875 -> 2791.4075
$ws.Cells.Item(122, 9).Value = 2760.3809  # I122: 2071.3333 -> 2760.3809
$ws.Cells.Item(122, 10).Value = 2900  # J122: 2137.5 -> 2900
$ws.Cells.Item(122, 11).Value = 8281.1427  # K122: 6213.999899999999 -> 8281.1427
$ws.Cells.Item(122, 12).Value = 8700  # L122: 6412.5 -> 8700
$ws.Cells.Item(122, 13).Value = -5831.1427  # M122: -3763.999899999999 -> -5831.1427
$ws.Cells.Item(122, 14).Value = -13600  # N122: -11312.5 -> -13600

# ALC!row127
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(127, 8).Value = 901.5  # H127: 879.4 -> 901.5
$ws.Cells.Item(127, 9).Value = 949.6667  # I127: 932.8333 -> 949.6667
$ws.Cells.Item(127, 10).Value = 757  # J127: 799.25 -> 757
$ws.Cells.Item(127, 11).Value = 2849.0001  # K127: 2798.4999 -> 2849.0001
$ws.Cells.Item(127, 12).Value = 2271  # L127: 2397.75 -> 2271
$ws.Cells.Item(127, 13).Value = 2110.9999  # M127: 2161.5001 -> 2110.9999
$ws.Cells.Item(127, 14).Value = -12191  # N127: -12317.75 -> -12191

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2783.5  # H137: 2889.4575 -> 2783.5
$ws.Cells.Item(137, 9).Value = 2858.4092  # I137: 3016.366 -> 2858.4092
$ws.Cells.Item(137, 11).Value = 8575.2276  # K137: 9049.098 -> 8575.2276
$ws.Cells.Item(137, 13).Value = -6025.2276  # M137: -6499.098 -> -6025.2276

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2353.724  # H138: 1803.5177 -> 2353.724
$ws.Cells.Item(138, 9).Value = 1811.35  # I138: 1103.0714 -> 1811.35
$ws.Cells.Item(138, 10).Value = 2639.1843  # J138: 2487.6743 -> 2639.1843
$ws.Cells.Item(138, 11).Value = 5434.049999999999  # K138: 3309.2142 -> 5434.049999999999
$ws.Cells.Item(138, 12).Value = 7917.5529  # L138: 7463.0229 -> 7917.5529
$ws.Cells.Item(138, 13).Value = -294.0499999999993  # M138: 1830.7858 -> -294.0499999999993
$ws.Cells.Item(138, 14).Value = -18197.5529  # N138: -17743.0229 -> -18197.5529

# ALC!row141
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 891765.5600000001  # H141: 836092.7 -> 891765.5600000001
$ws.Cells.Item(141, 9).Value = 2722.625  # I141: 2531.2222 -> 2722.625
$ws.Cells.Item(141, 11).Value = 8167.875  # K141: 7593.6666 -> 8167.875
$ws.Cells.Item(141, 13).Value = -2987.875  # M141: -2413.6666 -> -2987.875

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1458.5  # H74: 1801.8334 -> 1458.5
$ws.Cells.Item(74, 9).Value = 1361.1428  # I74: 1600.8889 -> 1361.1428
$ws.Cells.Item(74, 10).Value = 1685.6666  # J74: 2404.6667 -> 1685.6666
$ws.Cells.Item(74, 11).Value = 1361.1428  # K74: 1600.8889 -> 1361.1428
$ws.Cells.Item(74, 12).Value = 1685.6666  # L74: 2404.6667 -> 1685.6666
$ws.Cells.Item(74, 13).Value = -487.1428000000001  # M74: -726.8888999999999 -> -487.1428000000001
$ws.Cells.Item(74, 14).Value = -3433.6666  # N74: -4152.6667 -> -3433.6666

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 1458.5  # H77: 1801.8334 -> 1458.5
$ws.Cells.Item(77, 9).Value = 1361.1428  # I77: 1600.8889 -> 1361.1428
$ws.Cells.Item(77, 10).Value = 1685.6666  # J77: 2404.6667 -> 1685.6666
$ws.Cells.Item(77, 11).Value = 6805.714  # K77: 8004.4445 -> 6805.714
$ws.Cells.Item(77, 12).Value = 8428.333000000001  # L77: 12023.3335 -> 8428.333000000001
$ws.Cells.Item(77, 13).Value = -2437.714  # M77: -3636.4445 -> -2437.714
$ws.Cells.Item(77, 14).Value = -17164.333  # N77: -20759.3335 -> -17164.333

# BSM!row10
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(10, 8).Value = 60504.5  # H10: 73893.25 -> 60504.5
$ws.Cells.Item(10, 9).Value = 2000  # I10: 0 -> 2000
$ws.Cells.Item(10, 10).Value = 80006  # J10: 73893.25 -> 80006
$ws.Cells.Item(10, 11).Value = 2000  # K10: 0 -> 2000
$ws.Cells.Item(10, 12).Value = 80006  # L10: 73893.25 -> 80006
$ws.Cells.Item(10, 13).Value = -1860  # M10: None -> -1860
$ws.Cells.Item(10, 14).Value = -80286  # N10: -74173.25 -> -80286

# CRP!row9
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(9, 8).Value = 0  # H9: 1100 -> 0
$ws.Cells.Item(9, 10).Value = 0  # J9: 1100 -> 0
$ws.Cells.Item(9, 12).Value = 0  # L9: 1100 -> 0
$ws.Cells.Item(9, 14).ClearContents()  # N9: -1436 -> (removed)

# CRP!row86
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 14000  # H86: 8669.166999999999 -> 14000
$ws.Cells.Item(86, 9).Value = 13500  # I86: 10333.333 -> 13500
$ws.Cells.Item(86, 10).Value = 15000  # J86: 7005 -> 15000
$ws.Cells.Item(86, 11).Value = 13500  # K86: 10333.333 -> 13500
$ws.Cells.Item(86, 12).Value = 15000  # L86: 7005 -> 15000
$ws.Cells.Item(86, 13).Value = -12377  # M86: -9210.333000000001 -> -12377
$ws.Cells.Item(86, 14).Value = -17246  # N86: -9251 -> -17246

# CRP!row89
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 14000  # H89: 8669.166999999999 -> 14000
$ws.Cells.Item(89, 9).Value = 13500  # I89: 10333.333 -> 13500
$ws.Cells.Item(89, 10).Value = 15000  # J89: 7005 -> 15000
$ws.Cells.Item(89, 11).Value = 67500  # K89: 51666.665 -> 67500
$ws.Cells.Item(89, 12).Value = 75000  # L89: 35025 -> 75000
$ws.Cells.Item(89, 13).Value = -61884  # M89: -46050.665 -> -61884
$ws.Cells.Item(89, 14).Value = -86232  # N89: -46257 -> -86232

# CRP!row94
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 5873.8823  # H94: 5292.8 -> 5873.8823
$ws.Cells.Item(94, 10).Value = 5287  # J94: 4629.6 -> 5287
$ws.Cells.Item(94, 12).Value = 5287  # L94: 4629.6 -> 5287
$ws.Cells.Item(94, 14).Value = -6189  # N94: -5531.6 -> -6189

# CRP!row105
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 3042.1667  # H105: 4144.2856 -> 3042.1667
$ws.Cells.Item(105, 9).Value = 2440  # I105: 3502.5 -> 2440
$ws.Cells.Item(105, 10).Value = 4246.5  # J105: 5000 -> 4246.5
$ws.Cells.Item(105, 11).Value = 2440  # K105: 3502.5 -> 2440
$ws.Cells.Item(105, 12).Value = 4246.5  # L105: 5000 -> 4246.5
$ws.Cells.Item(105, 13).Value = -693  # M105: -1755.5 -> -693
$ws.Cells.Item(105, 14).Value = -7740.5  # N105: -8494 -> -7740.5

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 3112.24  # H132: 3564.2354 -> 3112.24
$ws.Cells.Item(132, 9).Value = 2639.4666  # I132: 2949.2 -> 2639.4666
$ws.Cells.Item(132, 10).Value = 3821.4  # J132: 4442.857 -> 3821.4
$ws.Cells.Item(132, 11).Value = 7918.399800000001  # K132: 8847.599999999999 -> 7918.399800000001
$ws.Cells.Item(132, 12).Value = 11464.2  # L132: 13328.571 -> 11464.2
$ws.Cells.Item(132, 13).Value = -5388.399800000001  # M132: -6317.599999999999 -> -5388.399800000001
$ws.Cells.Item(132, 14).Value = -16524.2  # N132: -18388.571 -> -16524.2

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 2794.1177  # H134: 4162.7 -> 2794.1177
$ws.Cells.Item(134, 9).Value = 886.5  # I134: 899.5 -> 886.5
$ws.Cells.Item(134, 10).Value = 4489.778  # J134: 4978.5 -> 4489.778
$ws.Cells.Item(134, 11).Value = 2659.5  # K134: 2698.5 -> 2659.5
$ws.Cells.Item(134, 12).Value = 13469.334  # L134: 14935.5 -> 13469.334
$ws.Cells.Item(134, 13).Value = -124.5  # M134: -163.5 -> -124.5
$ws.Cells.Item(134, 14).Value = -18539.334  # N134: -20005.5 -> -18539.334

# CUL!row34
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 6435.2383  # H34: 9598.333000000001 -> 6435.2383
$ws.Cells.Item(34, 9).Value = 425  # I34: 383.33334 -> 425
$ws.Cells.Item(34, 10).Value = 7849.4116  # J34: 18813.334 -> 7849.4116
$ws.Cells.Item(34, 11).Value = 1275  # K34: 1150.00002 -> 1275
$ws.Cells.Item(34, 12).Value = 23548.2348  # L34: 56440.00199999999 -> 23548.2348
$ws.Cells.Item(34, 13).Value = -1191  # M34: -1066.00002 -> -1191
$ws.Cells.Item(34, 14).Value = -23716.2348  # N34: -56608.00199999999 -> -23716.2348

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4011.9285  # H132: 4052.7666 -> 4011.9285
$ws.Cells.Item(132, 9).Value = 5183.6  # I132: 4363.2856 -> 5183.6
$ws.Cells.Item(132, 10).Value = 3361  # J132: 3781.0625 -> 3361
$ws.Cells.Item(132, 11).Value = 15550.8  # K132: 13089.8568 -> 15550.8
$ws.Cells.Item(132, 12).Value = 10083  # L132: 11343.1875 -> 10083
$ws.Cells.Item(132, 13).Value = -13020.8  # M132: -10559.8568 -> -13020.8
$ws.Cells.Item(132, 14).Value = -15143  # N132: -16403.1875 -> -15143

# LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1301.7413  # H46: 1411.5385 -> 1301.7413
$ws.Cells.Item(46, 9).Value = 906.02  # I46: 979.5454999999999 -> 906.02
$ws.Cells.Item(46, 10).Value = 3775  # J46: 3787.5 -> 3775
$ws.Cells.Item(46, 11).Value = 906.02  # K46: 979.5454999999999 -> 906.02
$ws.Cells.Item(46, 12).Value = 3775  # L46: 3787.5 -> 3775
$ws.Cells.Item(46, 13).Value = -718.02  # M46: -791.5454999999999 -> -718.02
$ws.Cells.Item(46, 14).Value = -4151  # N46: -4163.5 -> -4151

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2560.853  # H132: 3554.762 -> 2560.853
$ws.Cells.Item(132, 9).Value = 1642.76  # I132: 2387.5 -> 1642.76
$ws.Cells.Item(132, 11).Value = 4928.28  # K132: 7162.5 -> 4928.28
$ws.Cells.Item(132, 13).Value = -2398.28  # M132: -4632.5 -> -2398.28

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 2353.2104  # H136: 2189 -> 2353.2104
$ws.Cells.Item(136, 9).Value = 1744.88  # I136: 1597.4828 -> 1744.88
$ws.Cells.Item(136, 10).Value = 3523.077  # J136: 3414.2856 -> 3523.077
$ws.Cells.Item(136, 11).Value = 5234.64  # K136: 4792.4484 -> 5234.64
$ws.Cells.Item(136, 12).Value = 10569.231  # L136: 10242.8568 -> 10569.231
$ws.Cells.Item(136, 13).Value = -2684.64  # M136: -2242.4484 -> -2684.64
$ws.Cells.Item(136, 14).Value = -15669.231  # N136: -15342.8568 -> -15669.231

# WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 3848476  # H126: 3335642.8 -> 3848476
$ws.Cells.Item(126, 9).Value = 1693.7  # I126: 1626.091 -> 1693.7
$ws.Cells.Item(126, 10).Value = 16671084  # J126: 12504189 -> 16671084
$ws.Cells.Item(126, 11).Value = 5081.1  # K126: 4878.272999999999 -> 5081.1
$ws.Cells.Item(126, 12).Value = 50013252  # L126: 37512567 -> 50013252
$ws.Cells.Item(126, 13).Value = -2611.1  # M126: -2408.272999999999 -> -2611.1
$ws.Cells.Item(126, 14).Value = -50018192  # N126: -37517507 -> -50018192

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3628.362  # H132: 4131.04 -> 3628.362
$ws.Cells.Item(132, 9).Value = 1498.7916  # I132: 1679.738 -> 1498.7916
$ws.Cells.Item(132, 10).Value = 13850.3  # J132: 17000.375 -> 13850.3
$ws.Cells.Item(132, 11).Value = 4496.3748  # K132: 5039.214 -> 4496.3748
$ws.Cells.Item(132, 12).Value = 41550.89999999999  # L132: 51001.125 -> 41550.89999999999
$ws.Cells.Item(132, 13).Value = -1966.3748  # M132: -2509.214 -> -1966.3748
$ws.Cells.Item(132, 14).Value = -46610.89999999999  # N132: -56061.125 -> -46610.89999999999

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1118.826  # H136: 1141.4667 -> 1118.826
$ws.Cells.Item(136, 10).Value = 1905.8823  # J136: 2018.75 -> 1905.8823
$ws.Cells.Item(136, 12).Value = 5717.6469  # L136: 6056.25 -> 5717.6469
$ws.Cells.Item(136, 14).Value = -10817.6469  # N136: -11156.25 -> -10817.6469

# WVR!row140
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(140, 8).Value = 40769.89  # H140: 38344.273 -> 40769.89
$ws.Cells.Item(140, 10).Value = 40769.89  # J140: 38344.273 -> 40769.89
$ws.Cells.Item(140, 12).Value = 40769.89  # L140: 38344.273 -> 40769.89
$ws.Cells.Item(140, 14).Value = -51129.89  # N140: -48704.273 -> -51129.89
